# Refresh the "Price" (column D) quotes for the crypto symbol list.
# These cells hold their numbers as literal text (e.g. "241.41"), so a
# plain numeric assignment would let Excel re-interpret them as real
# numbers and silently mangle the formatting (trailing zeros dropped,
# floating-point noise appended, etc). Prefixing the value with a
# leading apostrophe forces Excel to store it verbatim as text, exactly
# like the source feed does; resetting the cell Style back to "Normal"
# afterwards clears the quote-prefix flag Excel would otherwise leave
# behind so the cell's formatting is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    "D2"  = "241.30";
    "D3"  = "21.74";
    "D4"  = "5.360";
    "D5"  = "0.05667";
    "D6"  = "3.411";
    "D8"  = "0.8068";
    "D9"  = "0.8579";
    "D10" = "0.1422";
    "D11" = "0.07255";
    "D12" = "0.03060";
    "D13" = "0.03155";
    "D14" = "0.09351";
    "D15" = "3.905";
    "D16" = "0.001585";
    "D17" = "0.04814";
    "D18" = "0.0005829";
    "D19" = "0.006310";
    "D20" = "0.0009964";
    "D22" = "0.0001502";
    "D23" = "3.740";
    "D24" = "2.148";
    "D26" = "0.1312";
    "D27" = "0.0004005";
    "D40" = "0.03792";
    "D41" = "0.006707";
    "D42" = "0.1046";
    "D43" = "0.002825";
    "D44" = "0.006829";
    "D45" = "0.00005621";
    "D47" = "0.5809";
    "D49" = "0.00002103";
    "D50" = "0.01012"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $priceUpdates[$addr]
    $cell.Style = "Normal"
}
